# M09 Froze Encoder 12345
# Updates Epoch Accuracy column B (accuracy values) for the re-run epochs,
# and moves the active cell to A2 within the existing full-sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 105, 106, 107, 108, 109, 110, 111, 112, 114, 115, 116, 117, 118)
$vals = @(0.578125, 0.484375, 0.46875, 0.390625, 0.421875, 0.390625, 0.265625, 0.34375, 0.265625, 0.265625, 0.296875, 0.328125, 0.265625, 0.265625, 0.25, 0.25, 0.234375, 0.21875, 0.21875, 0.203125, 0.203125, 0.203125, 0.21875, 0.203125, 0.203125, 0.25, 0.25, 0.25, 0.25, 0.234375, 0.234375, 0.234375, 0.25, 0.25, 0.25, 0.25, 0.25, 0.234375, 0.234375, 0.234375, 0.234375, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.21875, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.203125, 0.234375, 0.21875, 0.25, 0.21875, 0.203125, 0.15625, 0.234375, 0.203125, 0.21875, 0.1875, 0.3125, 0.1967213114754098)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $ws.Cells.Item($rows[$i], 2).Value = $vals[$i]
}

# Restore the whole-sheet selection with the active cell at A2.
$ws.Range("A1:XFD1048576").Select()
